$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A48").Value = @'
45.You are creating an online chat room using HTML WebSocket API.You want to receive WebSocket messages.Which event type should you listen to? 
'@

$ws.Range("B48").Value = @'
message
'@

$ws.Range("D48").Value = @'
My answer is correct. Answer explanation: You should use the message event.Event message occurs when  the client receives data from the web socket.
'@

$ws.Range("A49").Value = @'
46.You need to use jQuery ajax function to invoke the web service for exam number 70-480 and display the price as an alert message.
'@

$ws.Range("B49").Value = @'
$.ajax({
    url: "http://services.measureup.com/Exam",
    type: "GET",
    data: {examNumber: "70-480"},
    contentType: "application/json; charset=utf-8",
    dataType: "json",
    success: function (result) {
        window.alert(result.price);
    }
});
'@

$ws.Range("D49").Value = @'
My answer is correct. Answer explanation: The ajax function accepts a set of key-value pairs that specify the settings for an Asynchronous JavaScript and XML(AJAX) request.The dataType setting specifies the type of data returned from the web service.In this scenario the data type returned is JavaScript Object Notation(JSON).The success setting specifies a callback function to execute if the request is successful.The first parameter of the callback function is an object that represents the data returned.The object is formatted according to the type specified by the dataType setting.
'@

$ws.Range("B50").Value = @'
var socket = new WebSocket("wss://ChatService");
socket.onmessage = function (event) {
    document.writeln(event.data);
};
'@

$ws.Range("C50").Value = @'
var socket = new WebSocket("wss://ChatService");
socket.send("{document.writeln(socket.bufferedAmount)}");
'@

$ws.Range("A50").Value = @'
47.You need to implement a callback to write a message to the browser when it is received from the server.
'@

$ws.Range("D50").Value = @'
My answer is incorrect. Answer explanation: This code uses an anonymous function as a callback to handle the onmessage event of the WebSocket object.This event is raised when a message is received.The callback function accepts an event object as a parameter.The data property of the event object represents the message that is received.
'@

$ws.Range("A51").Value = @'
48.You need to call the GetCurrentLocation and implement callbacks to display the position if the position is found or to display the error message if the position is not found.
'@

$ws.Range("B51").Value = @'
GetCurrentLocation(
    function (position) {
        alert(position);
    },
    function (error) {
        alert(error);
    }
);
'@

$ws.Range("C51").Value = @'
function onSuccess(position) {
    alert(position);
}
function onError(error) {
    alert(error);
}
GetCurrentLocation(onSuccess(), onError());
'@

$ws.Range("D51").Value = @'
My answer is incorrect.Answer explanation: The GetCurrentLocation function accepts two parameters.The first is a reference to a callback function that is called if a position is successfully obtained.The second is a reference to a callback function that is called if a position is not obtained.The signature of each function takes one parameter as indicated in the body of the GetCurrentLocation function.This code uses anonymous functions as the callback functions.
'@

$ws.Range("A52").Value = @'
49.You need to modify the markup so that the alert message correctly displays the value of the speed variable.
'@

$ws.Range("B52").Value = @'
1. Replace the line 11 with this code:  window.alert(newSpeed);   2. Add this code between lines 08 and 09: var newSpeed = this.speed;                                                                                  
'@

$ws.Range("C52").Value = @'
Add this code between lines 01 and 02: var newSpeed = 0;
'@

$ws.Range("D52").Value = @'
My answer is incorrect.Answer explanation: 1.This code defines a new local variable named newSpeed that is equal to the speed variable associated with the function. 2. This code uses the newSpeed variable to display the message.This is necessary because the previous code uses the "this" keyword to refer to the speed variable.However the speed variable is not available in the context of the anonymous function because it was declared in a different scope.To make a variable available to an anonymous function you should create a new local variable.
'@

$ws.Range("A53").Value = @'
50.You need to call the Hide function to an object named button and display a message box after the object is hidden.
'@

$ws.Range("B53").Value = @'
Hide(button, function() {alert("Hidden");});
'@

$ws.Range("D53").Value = @'
My answer is correct.Answer explanation: The Hide function accepts two parameters.The first parameter is an object reference of the element to be hidden.The second parameter is the callback.If you examine the signature of the callback it is called with no parameters.Therefore you should pass an annonymous function with zero parameters as the second parameter to the Hide function.
'@

$ws.Range("A54").Value = 51

$ws.Range("A54").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 52
$win.ScrollColumn = 1

Write-Host "Done"
